$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 206 (shifts existing rows 206:304 down to 207:305)
$ws.Rows(206).Insert()

# Populate the newly inserted row with the new record
$ws.Range("A206").Value = 3
$ws.Range("B206").Value = "Femacal de La Calera"
$ws.Range("C206").Value = "Coquimbo"
$ws.Range("D206").Value = 44510
$ws.Range("E206").Value = 5
$ws.Range("F206").Value = 100112021
$ws.Range("G206").Value = "Ají"
$ws.Range("H206").Value = "Americana (o)"
$ws.Range("I206").Value = "Primera"
$ws.Range("J206").Value = 70
$ws.Range("K206").Value = 35000
$ws.Range("L206").Value = 36000
$ws.Range("M206").Value = 35500
$ws.Range("N206").Value = "`$/caja 15 kilos"
$ws.Range("O206").Value = "Limache"
$ws.Range("P206").Value = 2367
$ws.Range("Q206").Value = 15
$ws.Range("R206").Value = "Hortaliza"

# Make sure the date cell keeps the same date/time number format used by the rest of column D
$ws.Range("D206").NumberFormat = "YYYY-MM-DD HH:MM:SS"
